# Change Track.xlsx - "Update for initial Project Name Practice 4"
# Add a new change-log entry (row 13) to the Tabelle1 worksheet, matching
# the formatting of the preceding rows, and move the active selection
# past the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Copy the formatting of the last existing data row (row 12) down onto the
# new row 13 so the date/style/alignment match the rest of the table.
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)

# Fill in the new change-log entry.
$ws.Range("A13").Value = 41755          # 4/26/2014
$ws.Range("B13").Value = "11"
$ws.Range("C13").Value = "JEB"
$ws.Range("D13").Value = "Initial Base Project "
$ws.Range("E13").Value = "Done"

# Move the selection to the next empty row, as after data entry.
$ws.Range("A14").Select() | Out-Null
